$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing columns A,B,C,D (Nombre
# Completo, Documentos, Alias, Accion) shift right to B,C,D,E.
$ws.Columns.Item(1).Insert()

# The new "ID OFAC" column holds numeric-looking identifiers that must be
# stored as text, so mark the data rows as Text before writing the values.
$ws.Range("A2:A7").NumberFormat = "@"

# Header cell, styled the same as the other header cells on row 1.
$ws.Cells.Item(1, 1).Value = "ID OFAC"
$ws.Cells.Item(1, 2).Copy()
$ws.Cells.Item(1, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New OFAC ids for each existing row.
$ids = @(52752, 52753, 52754, 52755, 52756, 52757)
for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = [string]$ids[$i]
}

# Match the target column width for the new column (character width 9).
$ws.Columns.Item(1).ColumnWidth = 8.14
